$d = $word.ActiveDocument

# Heading text: drop the leading description, keep just "-esttab-"
$d.Content.Find.Execute("Produce a table from saved dataset from -esttab-", $true, $false, $false, $false, $false, $true, 1, $false, "-esttab-", 2)

# Translate recurring table header/label/value cells (applies to both tables,
# since both tables share the exact same labels and values).
$d.Content.Find.Execute("Model 1", $true, $false, $false, $false, $false, $true, 1, $false, "模型 1", 2)
$d.Content.Find.Execute("Model 2", $true, $false, $false, $false, $false, $true, 1, $false, "模型 2", 2)

$d.Content.Find.Execute("Weight (lbs.)", $true, $false, $false, $false, $false, $true, 1, $false, "重量(公斤)", 2)
$d.Content.Find.Execute("0.001***", $true, $false, $false, $false, $false, $true, 1, $false, "0.003***", 2)

$d.Content.Find.Execute("Gear Ratio", $true, $false, $false, $false, $false, $true, 1, $false, "变速比", 2)

$d.Content.Find.Execute("Turn Circle (ft.) ", $true, $false, $false, $false, $false, $true, 1, $false, "转弯半径(米) ", 2)
$d.Content.Find.Execute("0.024", $true, $false, $false, $false, $false, $true, 1, $false, "0.080", 2)
$d.Content.Find.Execute("0.061", $true, $false, $false, $false, $false, $true, 1, $false, "0.201", 2)

$d.Content.Find.Execute("Car type", $true, $false, $false, $false, $false, $true, 1, $false, "国籍", 2)
